# Apply updated "dSF" (column F) values for a set of rows, reflecting a
# repull/recalculation of the final-day stack delta data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = -5
    7  = -3
    8  = -4
    13 = 2
    16 = -8
    20 = -1
    30 = -2
    45 = 2
    51 = 1
    52 = -1
    60 = -2
    61 = 1
    68 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
